$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 3492.7144
$ws.Range("I18").Value = 3074.8333
$ws.Range("K18").Value = 3074.8333
$ws.Range("M18").Value = -2790.8333

$ws.Range("H33").Value = 136.15384
$ws.Range("I33").Value = 70
$ws.Range("K33").Value = 70
$ws.Range("M33").Value = 159

$ws.Range("H49").Value = 500
$ws.Range("I49").Value = 500
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 1500
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = -1364
$ws.Range("N49").Value = $null

$ws.Range("H113").Value = 4870
$ws.Range("I113").Value = 4870
$ws.Range("K113").Value = 4870
$ws.Range("M113").Value = -1616

$ws.Range("H132").Value = 1648.7273
$ws.Range("I132").Value = 1570.6666
$ws.Range("K132").Value = 4711.9998
$ws.Range("M132").Value = -2181.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7517.6577
$ws.Range("I32").Value = 5183.6562
$ws.Range("K32").Value = 5183.6562
$ws.Range("M32").Value = -4896.6562

$ws.Range("H45").Value = 2406.7856
$ws.Range("I45").Value = 2399.6155
$ws.Range("K45").Value = 2399.6155
$ws.Range("M45").Value = -2022.6155

$ws.Range("H76").Value = 18972.25
$ws.Range("J76").Value = 18972.25
$ws.Range("L76").Value = 18972.25
$ws.Range("N76").Value = -19648.25

$ws.Range("H79").Value = 18972.25
$ws.Range("J79").Value = 18972.25
$ws.Range("L79").Value = 18972.25
$ws.Range("N79").Value = -21312.25

$ws.Range("H92").Value = 63999
$ws.Range("J92").Value = 63999
$ws.Range("L92").Value = 63999
$ws.Range("N92").Value = -68991

$ws.Range("H97").Value = 1110.1428
$ws.Range("I97").Value = 1095.1666
$ws.Range("K97").Value = 1095.1666
$ws.Range("M97").Value = -599.1666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2653.647
$ws.Range("I105").Value = 2571.9167
$ws.Range("J105").Value = 2849.8
$ws.Range("K105").Value = 2571.9167
$ws.Range("L105").Value = 2849.8
$ws.Range("M105").Value = -824.9167000000002
$ws.Range("N105").Value = -6343.8

$ws.Range("H107").Value = 1537.7646
$ws.Range("I107").Value = 1508.4286
$ws.Range("K107").Value = 1508.4286
$ws.Range("M107").Value = 411.5714

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 11000
$ws.Range("I86").Value = 9000
$ws.Range("K86").Value = 9000
$ws.Range("M86").Value = -7877

$ws.Range("H89").Value = 11000
$ws.Range("I89").Value = 9000
$ws.Range("K89").Value = 45000
$ws.Range("M89").Value = -39384

$ws.Range("H99").Value = 12611.64
$ws.Range("I99").Value = 9727.1
$ws.Range("J99").Value = 14534.667
$ws.Range("K99").Value = 9727.1
$ws.Range("L99").Value = 14534.667
$ws.Range("M99").Value = -8229.1
$ws.Range("N99").Value = -17530.667

$ws.Range("H126").Value = 12611.64
$ws.Range("I126").Value = 9727.1
$ws.Range("J126").Value = 14534.667
$ws.Range("K126").Value = 29181.3
$ws.Range("L126").Value = 43604.001
$ws.Range("M126").Value = -26711.3
$ws.Range("N126").Value = -48544.001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 8471642
$ws.Range("I4").Value = 15732473
$ws.Range("J4").Value = 671.0833
$ws.Range("K4").Value = 47197419
$ws.Range("L4").Value = 2013.2499
$ws.Range("M4").Value = -47197307
$ws.Range("N4").Value = -2237.2499

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 12501.5
$ws.Range("J7").Value = 12501.5
$ws.Range("L7").Value = 12501.5
$ws.Range("N7").Value = -12725.5

$ws.Range("H8").Value = 12501.5
$ws.Range("J8").Value = 12501.5
$ws.Range("L8").Value = 12501.5
$ws.Range("N8").Value = -12779.5

$ws.Range("H21").Value = 11995
$ws.Range("I21").Value = 11995
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 11995
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -11822
$ws.Range("N21").Value = $null

$ws.Range("H30").Value = 11995
$ws.Range("I30").Value = 11995
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 11995
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -11890
$ws.Range("N30").Value = $null

$ws.Range("H36").Value = 9391.5
$ws.Range("I36").Value = 7550
$ws.Range("J36").Value = 10005.333
$ws.Range("K36").Value = 7550
$ws.Range("L36").Value = 10005.333
$ws.Range("M36").Value = -7065
$ws.Range("N36").Value = -10975.333

$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").Value = $null

$ws.Range("H46").Value = 16666.334

$ws.Range("H52").Value = 45000
$ws.Range("J52").Value = 45000
$ws.Range("L52").Value = 45000
$ws.Range("N52").Value = -45518

$ws.Range("H97").Value = 419.8
$ws.Range("I97").Value = 419.8
$ws.Range("K97").Value = 419.8
$ws.Range("M97").Value = 76.19999999999999

$ws.Range("H99").Value = 16235.5
$ws.Range("I99").Value = 2471
$ws.Range("J99").Value = 30000
$ws.Range("K99").Value = 2471
$ws.Range("L99").Value = 30000
$ws.Range("M99").Value = -225
$ws.Range("N99").Value = -34492

$ws.Range("H133").Value = 139979
$ws.Range("J133").Value = 139979
$ws.Range("L133").Value = 139979
$ws.Range("N133").Value = -150099

$ws.Range("H139").Value = 40000
$ws.Range("J139").Value = 40000
$ws.Range("L139").Value = 40000
$ws.Range("N139").Value = -50280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2205.2856
$ws.Range("I7").Value = 2205.2856
$ws.Range("K7").Value = 2205.2856
$ws.Range("M7").Value = -2093.2856

$ws.Range("H29").Value = 16000
$ws.Range("I29").Value = 16000
$ws.Range("K29").Value = 16000
$ws.Range("M29").Value = -15705

$ws.Range("H125").Value = 80000
$ws.Range("J125").Value = 80000
$ws.Range("L125").Value = 80000
$ws.Range("N125").Value = -89840

$ws.Range("H126").Value = 2205.2856
$ws.Range("I126").Value = 2205.2856
$ws.Range("K126").Value = 6615.8568
$ws.Range("M126").Value = -4145.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 1714.6666
$ws.Range("I3").Value = 150
$ws.Range("J3").Value = 2497
$ws.Range("K3").Value = 150
$ws.Range("L3").Value = 2497
$ws.Range("M3").Value = -36
$ws.Range("N3").Value = -2725

$ws.Range("H32").Value = 15979.8
$ws.Range("I32").Value = 15979.8
$ws.Range("K32").Value = 15979.8
$ws.Range("M32").Value = -15662.8

$ws.Range("H34").Value = 11000
$ws.Range("I34").Value = 11000
$ws.Range("K34").Value = 11000
$ws.Range("M34").Value = -10797

$ws.Range("H62").Value = 7000.077
$ws.Range("I62").Value = 5000
$ws.Range("J62").Value = 7889
$ws.Range("K62").Value = 5000
$ws.Range("L62").Value = 7889
$ws.Range("M62").Value = -4376
$ws.Range("N62").Value = -9137

$ws.Range("H65").Value = 7000.077
$ws.Range("I65").Value = 5000
$ws.Range("J65").Value = 7889
$ws.Range("K65").Value = 25000
$ws.Range("L65").Value = 39445
$ws.Range("M65").Value = -21880
$ws.Range("N65").Value = -45685

$ws.Range("H82").Value = 34999.5
$ws.Range("J82").Value = 34999.5
$ws.Range("L82").Value = 34999.5
$ws.Range("N82").Value = -35765.5

$ws.Range("H85").Value = 34999.5
$ws.Range("J85").Value = 34999.5
$ws.Range("L85").Value = 34999.5
$ws.Range("N85").Value = -37651.5

$ws.Range("H136").Value = 1908.8462
$ws.Range("I136").Value = 1035
$ws.Range("K136").Value = 3105
$ws.Range("M136").Value = -555
